# Coastal Surface Piercing Profilers - update ingest/cal sheet
# CP05MOAS-GL002 -> CP05MOAS-GL339, deployment number 3 -> 1
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # Moorings
$ws2 = $wb.Worksheets.Item(2)  # Asset_Cal_Info

# --- Asset_Cal_Info sheet ---
# ADCPAM block
$ws2.Range("A2").Value = "CP05MOAS-GL339-01-ADCPAM000"
$ws2.Range("A3").Value = "CP05MOAS-GL339-01-ADCPAM000"
$ws2.Range("A4").Value = "CP05MOAS-GL339-01-ADCPAM000"
$ws2.Range("A5").Value = "CP05MOAS-GL339-01-ADCPAM000"
$ws2.Range("C2:C5").Value = 1

# FLORTM block
$ws2.Range("A7").Value = "CP05MOAS-GL339-02-FLORTM000"
$ws2.Range("A8").Value = "CP05MOAS-GL339-02-FLORTM000"
$ws2.Range("A9").Value = "CP05MOAS-GL339-02-FLORTM000"
$ws2.Range("A10").Value = "CP05MOAS-GL339-02-FLORTM000"
$ws2.Range("C7:C10").Value = 1

# CTDGVM block
$ws2.Range("A12").Value = "CP05MOAS-GL339-03-CTDGVM000"
$ws2.Range("C12").Value = 1

# DOSTAM block
$ws2.Range("A14").Value = "CP05MOAS-GL339-04-DOSTAM000"
$ws2.Range("C14").Value = 1

# PARADM block
$ws2.Range("A16").Value = "CP05MOAS-GL339-05-PARADM000"
$ws2.Range("C16").Value = 1

# ENG000000 block
$ws2.Range("A18").Value = "CP05MOAS-GL339-00-ENG000000"
$ws2.Range("C18").Value = 1

# --- Moorings sheet ---
$ws1.Range("A2").Value = "CP05MOAS-GL339"
$ws1.Range("C2").Value = 1

# --- Restore view/selection state: active sheet switches to Moorings ---
$ws2.Range("C19").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C13").Select() | Out-Null
